$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1874.1111
$ws.Range("I43").Value = 1820.5
$ws.Range("J43").Value = 1917
$ws.Range("K43").Value = 1820.5
$ws.Range("L43").Value = 1917
$ws.Range("M43").Value = -1751.5
$ws.Range("N43").Value = -2055
$ws.Range("H76").Value = 9157.200000000001
$ws.Range("I76").Value = 8953.799999999999
$ws.Range("J76").Value = 9360.6
$ws.Range("K76").Value = 8953.799999999999
$ws.Range("L76").Value = 9360.6
$ws.Range("M76").Value = -8638.799999999999
$ws.Range("N76").Value = -9990.6
$ws.Range("H79").Value = 9157.200000000001
$ws.Range("I79").Value = 8953.799999999999
$ws.Range("J79").Value = 9360.6
$ws.Range("K79").Value = 8953.799999999999
$ws.Range("L79").Value = 9360.6
$ws.Range("M79").Value = -7861.799999999999
$ws.Range("N79").Value = -11544.6
$ws.Range("H80").Value = 953.8461
$ws.Range("I80").Value = 760
$ws.Range("K80").Value = 2280
$ws.Range("M80").Value = -1282
$ws.Range("H83").Value = 953.8461
$ws.Range("I83").Value = 760
$ws.Range("K83").Value = 6840
$ws.Range("M83").Value = -1848
$ws.Range("H86").Value = 1939
$ws.Range("I86").Value = 1922.7142
$ws.Range("K86").Value = 1922.7142
$ws.Range("M86").Value = -799.7141999999999
$ws.Range("H89").Value = 1939
$ws.Range("I89").Value = 1922.7142
$ws.Range("K89").Value = 9613.571
$ws.Range("M89").Value = -3997.571
$ws.Range("H112").Value = 69505.336
$ws.Range("J112").Value = 55801.895
$ws.Range("L112").Value = 167405.685
$ws.Range("N112").Value = -169621.685
$ws.Range("H129").Value = 2610.75
$ws.Range("J129").Value = 3124.5
$ws.Range("L129").Value = 9373.5
$ws.Range("N129").Value = -19373.5
$ws.Range("H132").Value = 1812.7435
$ws.Range("I132").Value = 1797.1389
$ws.Range("K132").Value = 5391.4167
$ws.Range("M132").Value = -2861.4167
$ws.Range("H137").Value = 5827.0586
$ws.Range("I137").Value = 5673.4443
$ws.Range("K137").Value = 17020.3329
$ws.Range("M137").Value = -14470.3329
$ws.Range("H138").Value = 3862.1025
$ws.Range("I138").Value = 2884
$ws.Range("K138").Value = 8652
$ws.Range("M138").Value = -3512

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 7015.857
$ws.Range("I16").Value = 3815.4
$ws.Range("K16").Value = 3815.4
$ws.Range("M16").Value = -3528.4
$ws.Range("H32").Value = 5073.552
$ws.Range("I32").Value = 2541.262
$ws.Range("J32").Value = 11720.8125
$ws.Range("K32").Value = 2541.262
$ws.Range("L32").Value = 11720.8125
$ws.Range("M32").Value = -2254.262
$ws.Range("N32").Value = -12294.8125
$ws.Range("H140").Value = 62870
$ws.Range("J140").Value = 62870
$ws.Range("L140").Value = 62870
$ws.Range("N140").Value = -73230

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 11572.556
$ws.Range("I7").Value = 25349.75
$ws.Range("J7").Value = 550.8
$ws.Range("K7").Value = 25349.75
$ws.Range("L7").Value = 550.8
$ws.Range("M7").Value = -25236.75
$ws.Range("N7").Value = -776.8
$ws.Range("H31").Value = 20492.715
$ws.Range("I31").Value = 21966
$ws.Range("K31").Value = 21966
$ws.Range("M31").Value = -21671
$ws.Range("H34").Value = 20492.715
$ws.Range("I34").Value = 21966
$ws.Range("K34").Value = 21966
$ws.Range("M34").Value = -21764
$ws.Range("H86").Value = 4543
$ws.Range("I86").Value = 4664
$ws.Range("J86").Value = 4422
$ws.Range("K86").Value = 4664
$ws.Range("L86").Value = 4422
$ws.Range("M86").Value = -3541
$ws.Range("N86").Value = -6668
$ws.Range("H89").Value = 4543
$ws.Range("I89").Value = 4664
$ws.Range("J89").Value = 4422
$ws.Range("K89").Value = 23320
$ws.Range("L89").Value = 22110
$ws.Range("M89").Value = -17704
$ws.Range("N89").Value = -33342
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 833.5
$ws.Range("I92").Value = 600.4
$ws.Range("J92").Value = 1999
$ws.Range("K92").Value = 1801.2
$ws.Range("L92").Value = 5997
$ws.Range("M92").Value = -553.1999999999998
$ws.Range("N92").Value = -8493
$ws.Range("H131").Value = 1767.5186
$ws.Range("I131").Value = 1014.86664
$ws.Range("K131").Value = 3044.59992
$ws.Range("M131").Value = 1995.40008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 412.22858
$ws.Range("I2").Value = 490.42307
$ws.Range("K2").Value = 490.42307
$ws.Range("M2").Value = -377.42307
$ws.Range("H107").Value = 673.15
$ws.Range("I107").Value = 483.46155
$ws.Range("K107").Value = 483.46155
$ws.Range("M107").Value = 1436.53845
$ws.Range("H138").Value = 110347.25
$ws.Range("J138").Value = 110333
$ws.Range("L138").Value = 110333
$ws.Range("N138").Value = -120613

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1549.2
$ws.Range("I19").Value = 682.3333
$ws.Range("J19").Value = 2849.5
$ws.Range("K19").Value = 682.3333
$ws.Range("L19").Value = 2849.5
$ws.Range("M19").Value = -512.3333
$ws.Range("N19").Value = -3189.5
$ws.Range("H55").Value = 390.5625
$ws.Range("I55").Value = 251.28572
$ws.Range("J55").Value = 498.8889
$ws.Range("K55").Value = 251.28572
$ws.Range("L55").Value = 498.8889
$ws.Range("M55").Value = -78.28572
$ws.Range("N55").Value = -844.8888999999999
$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680
$ws.Range("H132").Value = 26674532
$ws.Range("I132").Value = 26674532
$ws.Range("K132").Value = 80023596
$ws.Range("M132").Value = -80021066

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4267
$ws.Range("I14").Value = 1400.5
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 1400.5
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -1232.5
$ws.Range("N14").Value = -10336
$ws.Range("H22").Value = 15000
$ws.Range("J22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("N22").Value = -15586
$ws.Range("H32").Value = 14899.5
$ws.Range("I32").Value = 14899.5
$ws.Range("K32").Value = 14899.5
$ws.Range("M32").Value = -14582.5
$ws.Range("H62").Value = 14599.8
$ws.Range("I62").Value = 8999.5
$ws.Range("K62").Value = 8999.5
$ws.Range("M62").Value = -8375.5
$ws.Range("H65").Value = 14599.8
$ws.Range("I65").Value = 8999.5
$ws.Range("K65").Value = 44997.5
$ws.Range("M65").Value = -41877.5
$ws.Range("H126").Value = 2649.5
$ws.Range("I126").Value = 2649.5
$ws.Range("K126").Value = 7948.5
$ws.Range("M126").Value = -5478.5
$ws.Range("H132").Value = 20010408
$ws.Range("I132").Value = 20840008
$ws.Range("K132").Value = 62520024
$ws.Range("M132").Value = -62517494
